$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Taul1")

# Update classifier result value for E18
$ws.Range("E18").Value = 4.5

# Update computed ranking values on row 22
$ws.Range("E22").Value = 0.86667000000000005
$ws.Range("L22").Value = 0.87204000000000004

# Update ranking percentage for B28
$ws.Range("B28").Value = 0.7

# Update the active cell selection to reflect the new view state
$ws.Range("D15").Select()
